$d = $word.ActiveDocument

# Helper: insert a brand-new paragraph right after the paragraph at index
# $afterIndex, give it the text $text, and return the index of the newly
# created paragraph. Because InsertParagraphAfter() is called on a paragraph
# whose own mark already carries an explicit <w:jc w:val="left"/> (and, for
# the "anchor" paragraphs used below, an explicit <w:ind w:firstLine="708"/>
# too), the freshly minted paragraph mark inherits that same pPr verbatim -
# which is exactly the formatting the target paragraphs need.
function Insert-ParaAfter($afterIndex, $text) {
    $src = $d.Paragraphs($afterIndex)
    $src.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $d.Paragraphs($newIndex)
    $newPara.Range.Text = $text
    return $newIndex
}

# ---------------------------------------------------------------------
# Block C - goes right after the "... 19.1.26" paragraph (Paragraph 20),
# i.e. right before "Дедлайн 4:". Processed first (bottom-most block)
# so the paragraph indices used by the blocks above stay valid.
# ---------------------------------------------------------------------
$i = 20
$i = Insert-ParaAfter $i "Нарисованы и реализованы карты с системой врагов"
$i = Insert-ParaAfter $i "Система очков здоровья и система атаки "
$i = Insert-ParaAfter $i "Реализован сам игрок на показываемом экране"

# ---------------------------------------------------------------------
# Block B - goes right after the "... 28.12.25" paragraph (Paragraph 18),
# i.e. right before "Дедлайн 3:".
# ---------------------------------------------------------------------
$i = 18
$i = Insert-ParaAfter $i "Нарисованы основные действия главного персонажа"
$i = Insert-ParaAfter $i "Созданы стартовая локация и локация выбора персонажа/начала игры"
$i = Insert-ParaAfter $i "Создана локация настроек"
$i = Insert-ParaAfter $i "Подключены json файлы"

# ---------------------------------------------------------------------
# Block A - goes right after the "... 14.12.25" paragraph (Paragraph 16),
# i.e. right before "Дедлайн 2:".
# ---------------------------------------------------------------------
$i = 16
$i = Insert-ParaAfter $i "Организация всей команды и работы над проектом"
$i = Insert-ParaAfter $i "Начать реализовывать заготовки по персонажу и игре"
$i = Insert-ParaAfter $i "Подключить всю команду к проекту"

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
